$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would be auto-parsed as a number by Excel need a
# leading apostrophe to force text, then the quote-prefix style is reset
# back to Normal so no visible/format change is introduced.
function Set-TextValue($range, $value) {
    if ($value -match "^[+-]?[0-9]*\.?[0-9]+$") {
        $range.Value = "'" + $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}

# --- Row-by-row price / volume refresh ---
Set-TextValue $ws.Range('D2') '39.470.21'
Set-TextValue $ws.Range('E2') '  +2.02%  '
Set-TextValue $ws.Range('D3') '2.164.63'
Set-TextValue $ws.Range('E3') '  +3.81%  '
Set-TextValue $ws.Range('E4') '  +0.04%  '
Set-TextValue $ws.Range('D5') '229.02'
Set-TextValue $ws.Range('E5') '  +0.26%  '
Set-TextValue $ws.Range('D6') '0.623'
Set-TextValue $ws.Range('E6') '  +1.11%  '
Set-TextValue $ws.Range('D7') '63.49'
Set-TextValue $ws.Range('E7') '  +5.96%  '
Set-TextValue $ws.Range('E8') '  +0.00%  '
Set-TextValue $ws.Range('E9') '  +3.44%  '
Set-TextValue $ws.Range('D10') '0.0866'
Set-TextValue $ws.Range('E10') '  +2.75%  '
Set-TextValue $ws.Range('D11') '0.103'
Set-TextValue $ws.Range('E11') '  -0.27%  '
Set-TextValue $ws.Range('E12') '  +6.70%  '
Set-TextValue $ws.Range('D13') '2.485.24'
Set-TextValue $ws.Range('E13') '  +3.82%  '
Set-TextValue $ws.Range('E14') '  +1.59%  '
Set-TextValue $ws.Range('E15') '  +2.03%  '
Set-TextValue $ws.Range('E16') '  +1.67%  '
Set-TextValue $ws.Range('D17') '2.164.57'
Set-TextValue $ws.Range('E17') '  +3.87%  '
Set-TextValue $ws.Range('D18') '39.457.63'
Set-TextValue $ws.Range('E18') '  +2.13%  '
Set-TextValue $ws.Range('D19') '72.23'
Set-TextValue $ws.Range('E19') '  +1.05%  '
Set-TextValue $ws.Range('E20') '  +2.37%  '
Set-TextValue $ws.Range('D21') '0.0₃0854'
Set-TextValue $ws.Range('E21') '  +1.85%  '
Set-TextValue $ws.Range('D22') '229.22'
Set-TextValue $ws.Range('E22') '  +1.02%  '
Set-TextValue $ws.Range('E24') '  +1.51%  '
Set-TextValue $ws.Range('D25') '2.30'
Set-TextValue $ws.Range('E25') '  -3.14%  '
Set-TextValue $ws.Range('D26') '9.75'
Set-TextValue $ws.Range('E26') '  +2.33%  '
Set-TextValue $ws.Range('D27') '172.29'
Set-TextValue $ws.Range('E27') '  +0.78%  '
Set-TextValue $ws.Range('D28') '0.137'
Set-TextValue $ws.Range('E28') '  -0.79%  '
Set-TextValue $ws.Range('D29') '19.72'
Set-TextValue $ws.Range('E29') '  +3.02%  '
Set-TextValue $ws.Range('E30') '  -2.95%  '
Set-TextValue $ws.Range('E31') '  +10.26%  '
Set-TextValue $ws.Range('E32') '  +1.48%  '
Set-TextValue $ws.Range('E33') '  +3.86%  '
Set-TextValue $ws.Range('E36') '  +2.36%  '
Set-TextValue $ws.Range('E37') '  +1.81%  '
Set-TextValue $ws.Range('E38') '  +2.35%  '
Set-TextValue $ws.Range('E39') '  -0.01%  '
Set-TextValue $ws.Range('D40') '104.54'
Set-TextValue $ws.Range('E40') '  +4.27%  '
Set-TextValue $ws.Range('D41') '0.0231'
Set-TextValue $ws.Range('E41') '  +1.51%  '
Set-TextValue $ws.Range('E42') '  -0.66%  '
Set-TextValue $ws.Range('D43') '1.538.39'
Set-TextValue $ws.Range('E43') '  -0.40%  '
Set-TextValue $ws.Range('E44') '  +7.23%  '
Set-TextValue $ws.Range('D45') '0.0934'
Set-TextValue $ws.Range('E45') '  +1.01%  '
Set-TextValue $ws.Range('D49') '4.24'
Set-TextValue $ws.Range('E49') '  +3.40%  '
Set-TextValue $ws.Range('D50') '2.369.12'
Set-TextValue $ws.Range('E50') '  +3.80%  '
Set-TextValue $ws.Range('D51') '2.97'
Set-TextValue $ws.Range('E51') '  +0.33%  '

# --- Rows 46-48 rotate (ARBITRUM / HuobiToken / FraxShare reshuffled) ---
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D46') '7.92'
Set-TextValue $ws.Range('E46') '  +2.78%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D47') '1.11'
Set-TextValue $ws.Range('E47') '  +7.46%  '
$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D48') '2.81'
Set-TextValue $ws.Range('E48') '  -0.09%  '
